$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,3).Value = 41033
$ws.Cells.Item(2,4).Value = 59265984
$ws.Cells.Item(3,3).Value = 97573
$ws.Cells.Item(3,4).Value = 142949488
$ws.Cells.Item(4,3).Value = 33194
$ws.Cells.Item(4,4).Value = 49136467
$ws.Cells.Item(5,3).Value = 9446
$ws.Cells.Item(5,4).Value = 14031817
$ws.Cells.Item(6,3).Value = 2298
$ws.Cells.Item(6,4).Value = 3413920
$ws.Cells.Item(7,3).Value = 230
$ws.Cells.Item(7,4).Value = 340093
$ws.Cells.Item(12,3).Value = 44245
$ws.Cells.Item(12,4).Value = 59926705
$ws.Cells.Item(13,3).Value = 10386
$ws.Cells.Item(13,4).Value = 15003273
$ws.Cells.Item(14,3).Value = 27544
$ws.Cells.Item(14,4).Value = 40361723
$ws.Cells.Item(15,3).Value = 8754
$ws.Cells.Item(15,4).Value = 12992146
$ws.Cells.Item(16,3).Value = 2314
$ws.Cells.Item(16,4).Value = 3438653
$ws.Cells.Item(17,3).Value = 461
$ws.Cells.Item(17,4).Value = 680623
$ws.Cells.Item(18,3).Value = 44
$ws.Cells.Item(18,4).Value = 66000
$ws.Cells.Item(20,3).Value = 10858
$ws.Cells.Item(20,4).Value = 14315330
$ws.Cells.Item(21,3).Value = 14328
$ws.Cells.Item(21,4).Value = 20660112
$ws.Cells.Item(22,3).Value = 33418
$ws.Cells.Item(22,4).Value = 48998666
$ws.Cells.Item(23,3).Value = 10763
$ws.Cells.Item(23,4).Value = 15992729
$ws.Cells.Item(24,3).Value = 2845
$ws.Cells.Item(24,4).Value = 4228565
$ws.Cells.Item(25,3).Value = 583
$ws.Cells.Item(25,4).Value = 868592
$ws.Cells.Item(27,3).Value = 12382
$ws.Cells.Item(27,4).Value = 16476148
$ws.Cells.Item(28,3).Value = 8337
$ws.Cells.Item(28,4).Value = 12056951
$ws.Cells.Item(29,3).Value = 23984
$ws.Cells.Item(29,4).Value = 35188713
$ws.Cells.Item(30,3).Value = 8243
$ws.Cells.Item(30,4).Value = 12254512
$ws.Cells.Item(31,3).Value = 2101
$ws.Cells.Item(31,4).Value = 3134699
$ws.Cells.Item(32,3).Value = 416
$ws.Cells.Item(32,4).Value = 614749
$ws.Cells.Item(34,3).Value = 8889
$ws.Cells.Item(34,4).Value = 11723954
$ws.Cells.Item(35,3).Value = 3603
$ws.Cells.Item(35,4).Value = 5201957
$ws.Cells.Item(36,3).Value = 8499
$ws.Cells.Item(36,4).Value = 12416950
$ws.Cells.Item(37,3).Value = 3374
$ws.Cells.Item(37,4).Value = 5003238
$ws.Cells.Item(38,3).Value = 864
$ws.Cells.Item(38,4).Value = 1287055
$ws.Cells.Item(39,3).Value = 178
$ws.Cells.Item(39,4).Value = 264686
$ws.Cells.Item(41,3).Value = 2704
$ws.Cells.Item(41,4).Value = 3644830
$ws.Cells.Item(42,3).Value = 18586
$ws.Cells.Item(42,4).Value = 26839198
$ws.Cells.Item(43,3).Value = 54293
$ws.Cells.Item(43,4).Value = 79555674
$ws.Cells.Item(44,3).Value = 19917
$ws.Cells.Item(44,4).Value = 29569360
$ws.Cells.Item(45,3).Value = 5966
$ws.Cells.Item(45,4).Value = 8873462
$ws.Cells.Item(46,3).Value = 1377
$ws.Cells.Item(46,4).Value = 2055144
$ws.Cells.Item(47,3).Value = 81
$ws.Cells.Item(47,4).Value = 119515
$ws.Cells.Item(50,3).Value = 18001
$ws.Cells.Item(50,4).Value = 23860794
$ws.Cells.Item(51,3).Value = 2335
$ws.Cells.Item(51,4).Value = 3387531
$ws.Cells.Item(52,3).Value = 7862
$ws.Cells.Item(52,4).Value = 11544655
$ws.Cells.Item(53,3).Value = 2633
$ws.Cells.Item(53,4).Value = 3930733
$ws.Cells.Item(54,3).Value = 826
$ws.Cells.Item(54,4).Value = 1233915
$ws.Cells.Item(55,3).Value = 214
$ws.Cells.Item(55,4).Value = 316448
$ws.Cells.Item(56,3).Value = 22
$ws.Cells.Item(56,4).Value = 33000
$ws.Cells.Item(57,3).Value = 7854
$ws.Cells.Item(57,4).Value = 10809780
$ws.Cells.Item(58,3).Value = 1570
$ws.Cells.Item(58,4).Value = 3089456
$ws.Cells.Item(59,3).Value = 3798
$ws.Cells.Item(59,4).Value = 7472369
$ws.Cells.Item(60,3).Value = 1493
$ws.Cells.Item(60,4).Value = 2936260
$ws.Cells.Item(61,3).Value = 502
$ws.Cells.Item(61,4).Value = 987083
$ws.Cells.Item(62,3).Value = 178
$ws.Cells.Item(62,4).Value = 367100
$ws.Cells.Item(64,3).Value = 2440
$ws.Cells.Item(64,4).Value = 4452458
$ws.Cells.Item(65,3).Value = 16740
$ws.Cells.Item(65,4).Value = 24158988
$ws.Cells.Item(66,3).Value = 47792
$ws.Cells.Item(66,4).Value = 69866431
$ws.Cells.Item(67,3).Value = 16703
$ws.Cells.Item(67,4).Value = 24816037
$ws.Cells.Item(68,3).Value = 4863
$ws.Cells.Item(68,4).Value = 7243514
$ws.Cells.Item(69,3).Value = 1065
$ws.Cells.Item(69,4).Value = 1583199
$ws.Cells.Item(70,3).Value = 91
$ws.Cells.Item(70,4).Value = 133830
$ws.Cells.Item(73,3).Value = 15990
$ws.Cells.Item(73,4).Value = 21001383
$ws.Cells.Item(74,3).Value = 61020
$ws.Cells.Item(74,4).Value = 88714563
$ws.Cells.Item(75,3).Value = 166607
$ws.Cells.Item(75,4).Value = 245238903
$ws.Cells.Item(76,3).Value = 71007
$ws.Cells.Item(76,4).Value = 105761191
$ws.Cells.Item(77,3).Value = 23070
$ws.Cells.Item(77,4).Value = 34466443
$ws.Cells.Item(78,3).Value = 5806
$ws.Cells.Item(78,4).Value = 8669413
$ws.Cells.Item(79,3).Value = 406
$ws.Cells.Item(79,4).Value = 603972
$ws.Cells.Item(80,3).Value = 31
$ws.Cells.Item(80,4).Value = 45405
$ws.Cells.Item(85,3).Value = 60041
$ws.Cells.Item(85,4).Value = 81111781
$ws.Cells.Item(86,3).Value = 5075
$ws.Cells.Item(86,4).Value = 7352391
$ws.Cells.Item(87,3).Value = 12513
$ws.Cells.Item(87,4).Value = 18375992
$ws.Cells.Item(88,3).Value = 4106
$ws.Cells.Item(88,4).Value = 6118140
$ws.Cells.Item(89,3).Value = 1439
$ws.Cells.Item(89,4).Value = 2149611
$ws.Cells.Item(90,3).Value = 339
$ws.Cells.Item(90,4).Value = 505012
$ws.Cells.Item(91,3).Value = 31
$ws.Cells.Item(91,4).Value = 46402
$ws.Cells.Item(93,3).Value = 5833
$ws.Cells.Item(93,4).Value = 7819951
$ws.Cells.Item(94,3).Value = 1793
$ws.Cells.Item(94,4).Value = 2584083
$ws.Cells.Item(95,3).Value = 5743
$ws.Cells.Item(95,4).Value = 8464107
$ws.Cells.Item(96,3).Value = 2077
$ws.Cells.Item(96,4).Value = 3091281
$ws.Cells.Item(97,3).Value = 765
$ws.Cells.Item(97,4).Value = 1146460
$ws.Cells.Item(98,3).Value = 209
$ws.Cells.Item(98,4).Value = 315613
$ws.Cells.Item(101,3).Value = 3906
$ws.Cells.Item(101,4).Value = 5181853
$ws.Cells.Item(102,3).Value = 887
$ws.Cells.Item(102,4).Value = 1683061
$ws.Cells.Item(103,3).Value = 589
$ws.Cells.Item(103,4).Value = 1184992
$ws.Cells.Item(104,3).Value = 209
$ws.Cells.Item(104,4).Value = 408194
$ws.Cells.Item(105,3).Value = 68
$ws.Cells.Item(105,4).Value = 130500
$ws.Cells.Item(106,3).Value = 39
$ws.Cells.Item(106,4).Value = 84000
$ws.Cells.Item(107,3).Value = 11749
$ws.Cells.Item(107,4).Value = 17029004
$ws.Cells.Item(108,3).Value = 30875
$ws.Cells.Item(108,4).Value = 45325803
$ws.Cells.Item(109,3).Value = 10339
$ws.Cells.Item(109,4).Value = 15372249
$ws.Cells.Item(110,3).Value = 2865
$ws.Cells.Item(110,4).Value = 4270571
$ws.Cells.Item(111,3).Value = 549
$ws.Cells.Item(111,4).Value = 817453
$ws.Cells.Item(112,3).Value = 65
$ws.Cells.Item(112,4).Value = 97500
$ws.Cells.Item(114,3).Value = 3
$ws.Cells.Item(114,4).Value = 4500
$ws.Cells.Item(114,7).NumberFormat = "@"
$ws.Cells.Item(114,7).Value = "32"
$ws.Cells.Item(114,8).Value = "250 à 499 salariés"
$ws.Cells.Item(115,3).Value = 10401
$ws.Cells.Item(115,4).Value = 13692376
$ws.Cells.Item(115,5).NumberFormat = "@"
$ws.Cells.Item(115,5).Value = "28"
$ws.Cells.Item(115,6).Value = "Normandie"
$ws.Cells.Item(115,7).Value = "NN"
$ws.Cells.Item(115,8).Value = "Etablissement non employeur"
$ws.Cells.Item(116,3).Value = 32991
$ws.Cells.Item(116,4).Value = 47543917
$ws.Cells.Item(116,7).NumberFormat = "@"
$ws.Cells.Item(116,7).Value = "00"
$ws.Cells.Item(116,8).Value = "0 salarié"
$ws.Cells.Item(117,3).Value = 70260
$ws.Cells.Item(117,4).Value = 102769535
$ws.Cells.Item(117,7).NumberFormat = "@"
$ws.Cells.Item(117,7).Value = "01"
$ws.Cells.Item(117,8).Value = "1 ou 2 salariés"
$ws.Cells.Item(118,3).Value = 22530
$ws.Cells.Item(118,4).Value = 33464917
$ws.Cells.Item(118,7).NumberFormat = "@"
$ws.Cells.Item(118,7).Value = "02"
$ws.Cells.Item(118,8).Value = "3 à 5 salariés"
$ws.Cells.Item(119,3).Value = 6454
$ws.Cells.Item(119,4).Value = 9609534
$ws.Cells.Item(119,7).NumberFormat = "@"
$ws.Cells.Item(119,7).Value = "03"
$ws.Cells.Item(119,8).Value = "6 à 9 salariés"
$ws.Cells.Item(120,3).Value = 1269
$ws.Cells.Item(120,4).Value = 1895237
$ws.Cells.Item(120,7).NumberFormat = "@"
$ws.Cells.Item(120,7).Value = "11"
$ws.Cells.Item(120,8).Value = "10 à 19 salariés"
$ws.Cells.Item(121,3).Value = 115
$ws.Cells.Item(121,4).Value = 168895
$ws.Cells.Item(121,7).NumberFormat = "@"
$ws.Cells.Item(121,7).Value = "12"
$ws.Cells.Item(121,8).Value = "20 à 49 salariés"
$ws.Cells.Item(122,3).Value = 13
$ws.Cells.Item(122,4).Value = 19500
$ws.Cells.Item(122,7).NumberFormat = "@"
$ws.Cells.Item(122,7).Value = "21"
$ws.Cells.Item(122,8).Value = "50 à 99 salariés"
$ws.Cells.Item(123,3).Value = 5
$ws.Cells.Item(123,4).Value = 6070
$ws.Cells.Item(123,7).NumberFormat = "@"
$ws.Cells.Item(123,7).Value = "22"
$ws.Cells.Item(123,8).Value = "100 à 199 salariés"
$ws.Cells.Item(124,3).Value = 6
$ws.Cells.Item(124,4).Value = 9000
$ws.Cells.Item(124,7).NumberFormat = "@"
$ws.Cells.Item(124,7).Value = "32"
$ws.Cells.Item(124,8).Value = "250 à 499 salariés"
$ws.Cells.Item(125,3).Value = 27453
$ws.Cells.Item(125,4).Value = 36571849
$ws.Cells.Item(125,5).NumberFormat = "@"
$ws.Cells.Item(125,5).Value = "75"
$ws.Cells.Item(125,6).Value = "Nouvelle-Aquitaine"
$ws.Cells.Item(125,7).Value = "NN"
$ws.Cells.Item(125,8).Value = "Etablissement non employeur"
$ws.Cells.Item(126,3).Value = 39342
$ws.Cells.Item(126,4).Value = 56723234
$ws.Cells.Item(126,7).NumberFormat = "@"
$ws.Cells.Item(126,7).Value = "00"
$ws.Cells.Item(126,8).Value = "0 salarié"
$ws.Cells.Item(127,3).Value = 82403
$ws.Cells.Item(127,4).Value = 120419047
$ws.Cells.Item(127,7).NumberFormat = "@"
$ws.Cells.Item(127,7).Value = "01"
$ws.Cells.Item(127,8).Value = "1 ou 2 salariés"
$ws.Cells.Item(128,3).Value = 25292
$ws.Cells.Item(128,4).Value = 37536557
$ws.Cells.Item(128,7).NumberFormat = "@"
$ws.Cells.Item(128,7).Value = "02"
$ws.Cells.Item(128,8).Value = "3 à 5 salariés"
$ws.Cells.Item(129,3).Value = 6874
$ws.Cells.Item(129,4).Value = 10213650
$ws.Cells.Item(129,7).NumberFormat = "@"
$ws.Cells.Item(129,7).Value = "03"
$ws.Cells.Item(129,8).Value = "6 à 9 salariés"
$ws.Cells.Item(130,3).Value = 1440
$ws.Cells.Item(130,4).Value = 2136248
$ws.Cells.Item(130,7).NumberFormat = "@"
$ws.Cells.Item(130,7).Value = "11"
$ws.Cells.Item(130,8).Value = "10 à 19 salariés"
$ws.Cells.Item(131,3).Value = 82
$ws.Cells.Item(131,4).Value = 121228
$ws.Cells.Item(131,7).NumberFormat = "@"
$ws.Cells.Item(131,7).Value = "12"
$ws.Cells.Item(131,8).Value = "20 à 49 salariés"
$ws.Cells.Item(132,3).Value = 19
$ws.Cells.Item(132,4).Value = 28500
$ws.Cells.Item(132,7).NumberFormat = "@"
$ws.Cells.Item(132,7).Value = "21"
$ws.Cells.Item(132,8).Value = "50 à 99 salariés"
$ws.Cells.Item(133,3).Value = 3
$ws.Cells.Item(133,4).Value = 2100
$ws.Cells.Item(133,7).NumberFormat = "@"
$ws.Cells.Item(133,7).Value = "32"
$ws.Cells.Item(133,8).Value = "250 à 499 salariés"
$ws.Cells.Item(134,3).Value = 33879
$ws.Cells.Item(134,4).Value = 44881821
$ws.Cells.Item(134,5).NumberFormat = "@"
$ws.Cells.Item(134,5).Value = "76"
$ws.Cells.Item(134,6).Value = "Occitanie"
$ws.Cells.Item(134,7).Value = "NN"
$ws.Cells.Item(134,8).Value = "Etablissement non employeur"
$ws.Cells.Item(135,3).Value = 14353
$ws.Cells.Item(135,4).Value = 20766350
$ws.Cells.Item(135,7).NumberFormat = "@"
$ws.Cells.Item(135,7).Value = "00"
$ws.Cells.Item(135,8).Value = "0 salarié"
$ws.Cells.Item(136,3).Value = 34307
$ws.Cells.Item(136,4).Value = 50358677
$ws.Cells.Item(136,7).NumberFormat = "@"
$ws.Cells.Item(136,7).Value = "01"
$ws.Cells.Item(136,8).Value = "1 ou 2 salariés"
$ws.Cells.Item(137,3).Value = 12107
$ws.Cells.Item(137,4).Value = 17987073
$ws.Cells.Item(137,7).NumberFormat = "@"
$ws.Cells.Item(137,7).Value = "02"
$ws.Cells.Item(137,8).Value = "3 à 5 salariés"
$ws.Cells.Item(138,3).Value = 3201
$ws.Cells.Item(138,4).Value = 4771375
$ws.Cells.Item(138,7).NumberFormat = "@"
$ws.Cells.Item(138,7).Value = "03"
$ws.Cells.Item(138,8).Value = "6 à 9 salariés"
$ws.Cells.Item(139,3).Value = 565
$ws.Cells.Item(139,4).Value = 841490
$ws.Cells.Item(139,7).NumberFormat = "@"
$ws.Cells.Item(139,7).Value = "11"
$ws.Cells.Item(139,8).Value = "10 à 19 salariés"
$ws.Cells.Item(140,3).Value = 44
$ws.Cells.Item(140,4).Value = 64825
$ws.Cells.Item(140,7).NumberFormat = "@"
$ws.Cells.Item(140,7).Value = "12"
$ws.Cells.Item(140,8).Value = "20 à 49 salariés"
$ws.Cells.Item(141,3).Value = 8
$ws.Cells.Item(141,4).Value = 12000
$ws.Cells.Item(141,7).NumberFormat = "@"
$ws.Cells.Item(141,7).Value = "21"
$ws.Cells.Item(141,8).Value = "50 à 99 salariés"
$ws.Cells.Item(142,3).Value = 3
$ws.Cells.Item(142,4).Value = 984
$ws.Cells.Item(142,7).NumberFormat = "@"
$ws.Cells.Item(142,7).Value = "42"
$ws.Cells.Item(142,8).Value = "1 000 à 1 999 salariés"
$ws.Cells.Item(143,3).Value = 11493
$ws.Cells.Item(143,4).Value = 15282188
$ws.Cells.Item(143,5).NumberFormat = "@"
$ws.Cells.Item(143,5).Value = "52"
$ws.Cells.Item(143,6).Value = "Pays de la Loire"
$ws.Cells.Item(143,7).Value = "NN"
$ws.Cells.Item(143,8).Value = "Etablissement non employeur"
$ws.Cells.Item(144,3).Value = 38546
$ws.Cells.Item(144,4).Value = 55664432
$ws.Cells.Item(144,7).NumberFormat = "@"
$ws.Cells.Item(144,7).Value = "00"
$ws.Cells.Item(144,8).Value = "0 salarié"
$ws.Cells.Item(145,3).Value = 88273
$ws.Cells.Item(145,4).Value = 129247300
$ws.Cells.Item(145,7).NumberFormat = "@"
$ws.Cells.Item(145,7).Value = "01"
$ws.Cells.Item(145,8).Value = "1 ou 2 salariés"
$ws.Cells.Item(146,3).Value = 26267
$ws.Cells.Item(146,4).Value = 39009656
$ws.Cells.Item(146,7).NumberFormat = "@"
$ws.Cells.Item(146,7).Value = "02"
$ws.Cells.Item(146,8).Value = "3 à 5 salariés"
$ws.Cells.Item(147,3).Value = 6956
$ws.Cells.Item(147,4).Value = 10364072
$ws.Cells.Item(147,7).NumberFormat = "@"
$ws.Cells.Item(147,7).Value = "03"
$ws.Cells.Item(147,8).Value = "6 à 9 salariés"
$ws.Cells.Item(148,3).Value = 1634
$ws.Cells.Item(148,4).Value = 2427169
$ws.Cells.Item(148,7).NumberFormat = "@"
$ws.Cells.Item(148,7).Value = "11"
$ws.Cells.Item(148,8).Value = "10 à 19 salariés"
$ws.Cells.Item(149,3).Value = 107
$ws.Cells.Item(149,4).Value = 160130
$ws.Cells.Item(149,7).NumberFormat = "@"
$ws.Cells.Item(149,7).Value = "12"
$ws.Cells.Item(149,8).Value = "20 à 49 salariés"
$ws.Cells.Item(150,3).Value = 17
$ws.Cells.Item(150,4).Value = 25500
$ws.Cells.Item(150,7).NumberFormat = "@"
$ws.Cells.Item(150,7).Value = "21"
$ws.Cells.Item(150,8).Value = "50 à 99 salariés"
$ws.Cells.Item(151,1).Value = "Fonds de solidarité"
$ws.Cells.Item(151,2).Value = "VOLET1"
$ws.Cells.Item(151,3).Value = 31358
$ws.Cells.Item(151,4).Value = 42184940
$ws.Cells.Item(151,5).NumberFormat = "@"
$ws.Cells.Item(151,5).Value = "93"
$ws.Cells.Item(151,6).Value = "Provence-Alpes-Côte d'Azur"
$ws.Cells.Item(151,7).Value = "NN"
$ws.Cells.Item(151,8).Value = "Etablissement non employeur"
